$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price/volume/hour data per Jan 14 2023 21:07 UTC refresh.
# D/E/G columns are stored as text (inlineStr) in the source sheet, so we
# quote-prefix the values to stop Excel from re-typing them as numbers
# (which would strip meaningful trailing zeros like "3.820" -> 3.82).

$ws.Range("D2").Value = "'302.99"
$ws.Range("E2").Value = "'5.03%"
$ws.Range("G2").Value = "'21"

$ws.Range("D3").Value = "'32.15"
$ws.Range("E3").Value = "'9.92%"
$ws.Range("G3").Value = "'21"

$ws.Range("D4").Value = "'5.273"
$ws.Range("E4").Value = "'-0.21%"
$ws.Range("G4").Value = "'21"

$ws.Range("E5").Value = "'6.20%"
$ws.Range("G5").Value = "'21"

$ws.Range("D6").Value = "'7.875"
$ws.Range("E6").Value = "'5.72%"
$ws.Range("G6").Value = "'21"

$ws.Range("D7").Value = "'3.820"
$ws.Range("E7").Value = "'7.24%"
$ws.Range("G7").Value = "'21"

$ws.Range("D8").Value = "'1.509"
$ws.Range("E8").Value = "'8.16%"
$ws.Range("G8").Value = "'21"

$ws.Range("D9").Value = "'0.9193"
$ws.Range("E9").Value = "'1.36%"
$ws.Range("G9").Value = "'21"

$ws.Range("D10").Value = "'0.1693"
$ws.Range("E10").Value = "'5.27%"
$ws.Range("G10").Value = "'21"

$ws.Range("D11").Value = "'0.07855"
$ws.Range("E11").Value = "'2.63%"
$ws.Range("G11").Value = "'21"

$ws.Range("D12").Value = "'0.07999"
$ws.Range("E12").Value = "'3.85%"
$ws.Range("G12").Value = "'21"

$ws.Range("D13").Value = "'0.03073"
$ws.Range("E13").Value = "'5.83%"
$ws.Range("G13").Value = "'21"

$ws.Range("D14").Value = "'0.09888"
$ws.Range("E14").Value = "'9.52%"
$ws.Range("G14").Value = "'21"

$ws.Range("D15").Value = "'0.001486"
$ws.Range("E15").Value = "'-5.65%"
$ws.Range("G15").Value = "'21"

$ws.Range("D16").Value = "'0.04604"
$ws.Range("E16").Value = "'1.64%"
$ws.Range("G16").Value = "'21"

$ws.Range("D17").Value = "'0.006178"
$ws.Range("E17").Value = "'1.64%"
$ws.Range("G17").Value = "'21"

$ws.Range("D18").Value = "'3.460"
$ws.Range("E18").Value = "'-1.02%"
$ws.Range("G18").Value = "'21"

$ws.Range("D19").Value = "'2.231"
$ws.Range("E19").Value = "'-0.04%"
$ws.Range("G19").Value = "'21"

$ws.Range("D20").Value = "'0.3302"
$ws.Range("E20").Value = "'2.20%"
$ws.Range("G20").Value = "'21"

$ws.Range("D21").Value = "'0.1341"
$ws.Range("E21").Value = "'-0.50%"
$ws.Range("G21").Value = "'21"

$ws.Range("D22").Value = "'4.484"
$ws.Range("E22").Value = "'11.77%"
$ws.Range("G22").Value = "'21"

$ws.Range("D23").Value = "'0.1618"
$ws.Range("E23").Value = "'1.39%"
$ws.Range("G23").Value = "'21"

$ws.Range("D24").Value = "'0.001217"
$ws.Range("E24").Value = "'0.43%"
$ws.Range("G24").Value = "'21"

$ws.Range("E25").Value = "'6.77%"
$ws.Range("G25").Value = "'21"

$ws.Range("D26").Value = "'0.0001399"
$ws.Range("E26").Value = "'19.55%"
$ws.Range("G26").Value = "'21"

$ws.Range("D27").Value = "'0.0001777"
$ws.Range("E27").Value = "'6.72%"
$ws.Range("G27").Value = "'21"

$ws.Range("G28").Value = "'21"

$ws.Range("G29").Value = "'21"

$ws.Range("G30").Value = "'21"

$ws.Range("G31").Value = "'21"

$ws.Range("G32").Value = "'21"

$ws.Range("G33").Value = "'21"

$ws.Range("G34").Value = "'21"

$ws.Range("G35").Value = "'21"

$ws.Range("G36").Value = "'21"

$ws.Range("G37").Value = "'21"

$ws.Range("G38").Value = "'21"

$ws.Range("D39").Value = "'0.01726"
$ws.Range("E39").Value = "'2,550.35%"
$ws.Range("G39").Value = "'21"

$ws.Range("D40").Value = "'0.04486"
$ws.Range("E40").Value = "'2.54%"
$ws.Range("G40").Value = "'21"

$ws.Range("D41").Value = "'0.006894"
$ws.Range("E41").Value = "'-1.56%"
$ws.Range("G41").Value = "'21"

$ws.Range("E42").Value = "'7.73%"
$ws.Range("G42").Value = "'21"

$ws.Range("D43").Value = "'0.002199"
$ws.Range("E43").Value = "'6.19%"
$ws.Range("G43").Value = "'21"

$ws.Range("D44").Value = "'0.01286"
$ws.Range("E44").Value = "'9.57%"
$ws.Range("G44").Value = "'21"

$ws.Range("D45").Value = "'0.00006143"
$ws.Range("E45").Value = "'5.21%"
$ws.Range("G45").Value = "'21"

$ws.Range("D46").Value = "'0.7094"
$ws.Range("E46").Value = "'-63.23%"
$ws.Range("G46").Value = "'21"

$ws.Range("D47").Value = "'0.01496"
$ws.Range("E47").Value = "'15.40%"
$ws.Range("G47").Value = "'21"

$ws.Range("G48").Value = "'21"

$ws.Range("G49").Value = "'21"

$ws.Range("G50").Value = "'21"

$ws.Range("G51").Value = "'21"
